$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41 ("tce") gets an "x" in column B (CBW column) indicating that
# no part of that record will be processed, per the sheet's documented
# convention in row 6.
$ws.Range("B41").Value = "x"

# Move the active selection to the newly-edited cell (matches the
# author's view state after making the edit).
$ws.Range("B41").Select()
